# "removed false start data"
# The first two logged rows (for 2016-09-26 and 2016-09-27) were a false
# start and are removed; the remaining rows shift up so the sheet starts
# with the 2016-09-28 / 2016-09-29 entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two stale rows (worksheet rows 2 and 3) - this shifts the
# rows below them up, so former rows 4 and 5 become rows 2 and 3.
$ws.Rows("2:3").Delete()

# Leave the selection on the (new) last data row, matching the state the
# workbook was left in after the cleanup.
$ws.Range("A3").EntireRow.Select() | Out-Null
